$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7; existing rows 7-13 (weekly Alcachofa
# records) shift down to become rows 8-14, keeping their data intact.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(7, 3).Value = 'Ñuble'
$ws.Cells.Item(7, 4).Value = 44467
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112013
$ws.Cells.Item(7, 7).Value = 'Alcachofa'
$ws.Cells.Item(7, 8).Value = 'Madrigal'
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 160
$ws.Cells.Item(7, 11).Value = 11000
$ws.Cells.Item(7, 12).Value = 12000
$ws.Cells.Item(7, 13).Value = 11500
$ws.Cells.Item(7, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(7, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(7, 16).Value = 288
$ws.Cells.Item(7, 17).Value = 40
$ws.Cells.Item(7, 18).Value = 'Hortaliza'
